$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple numeric value updates ---
$ws.Range("L6").Value = 123
$ws.Range("L7").Value = 9

# --- Percentage-like text cell updates (preserve text type & style) ---
$ws.Range("L9").Value = "'38.0%"
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L10").Value = "'71.2%"
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("R21").Value = "'37.0%"
$ws.Range("K21").Copy()
$ws.Range("R21").PasteSpecial(-4122)
$ws.Range("S21").Value = "'76.5%"
$ws.Range("L21").Copy()
$ws.Range("S21").PasteSpecial(-4122)
$ws.Range("R22").Value = "'37.0%"
$ws.Range("K22").Copy()
$ws.Range("R22").PasteSpecial(-4122)
$ws.Range("S22").Value = "'73.7%"
$ws.Range("L22").Copy()
$ws.Range("S22").PasteSpecial(-4122)
$ws.Range("R23").Value = "'37.0%"
$ws.Range("K23").Copy()
$ws.Range("R23").PasteSpecial(-4122)
$ws.Range("S23").Value = "'80.0%"
$ws.Range("L23").Copy()
$ws.Range("S23").PasteSpecial(-4122)
$ws.Range("R24").Value = "'33.3%"
$ws.Range("K24").Copy()
$ws.Range("R24").PasteSpecial(-4122)
$ws.Range("S24").Value = "'67.5%"
$ws.Range("L24").Copy()
$ws.Range("S24").PasteSpecial(-4122)
$ws.Range("R25").Value = "'37.0%"
$ws.Range("K25").Copy()
$ws.Range("R25").PasteSpecial(-4122)
$ws.Range("S25").Value = "'66.5%"
$ws.Range("L25").Copy()
$ws.Range("S25").PasteSpecial(-4122)
$ws.Range("R26").Value = "'37.0%"
$ws.Range("K26").Copy()
$ws.Range("R26").PasteSpecial(-4122)
$ws.Range("S26").Value = "'59.0%"
$ws.Range("L26").Copy()
$ws.Range("S26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Class statistics O/P numeric updates ---
$ws.Range("O21").Value = 10
$ws.Range("P21").Value = 0
$ws.Range("O22").Value = 10
$ws.Range("P22").Value = 0
$ws.Range("O23").Value = 10
$ws.Range("P23").Value = 0
$ws.Range("O24").Value = 9
$ws.Range("P24").Value = 1
$ws.Range("O25").Value = 10
$ws.Range("P25").Value = 0
$ws.Range("O26").Value = 10
$ws.Range("P26").Value = 0

# --- Recorded-by text swaps ---
$ws.Range("G8").Value = "System, dnasr281@gmail.com"
$ws.Range("G9").Value = "System, dnasr281@gmail.com"
$ws.Range("G10").Value = "System, dnasr281@gmail.com"
$ws.Range("G35").Value = "System, dnasr281@gmail.com"
$ws.Range("G36").Value = "System, dnasr281@gmail.com"
$ws.Range("G37").Value = "System, dnasr281@gmail.com"
$ws.Range("G62").Value = "System, dnasr281@gmail.com"
$ws.Range("G63").Value = "System, dnasr281@gmail.com"
$ws.Range("G64").Value = "System, dnasr281@gmail.com"
$ws.Range("G89").Value = "System, dnasr281@gmail.com"
$ws.Range("G90").Value = "System, dnasr281@gmail.com"
$ws.Range("G91").Value = "System, dnasr281@gmail.com"
$ws.Range("G116").Value = "System, dnasr281@gmail.com"
$ws.Range("G117").Value = "System, dnasr281@gmail.com"
$ws.Range("G118").Value = "System, dnasr281@gmail.com"
$ws.Range("G143").Value = "System, dnasr281@gmail.com"
$ws.Range("G144").Value = "System, dnasr281@gmail.com"
$ws.Range("G145").Value = "System, dnasr281@gmail.com"
$ws.Range("G170").Value = "System, dnasr281@gmail.com"
$ws.Range("G197").Value = "System, dnasr281@gmail.com"
$ws.Range("G224").Value = "System, dnasr281@gmail.com"
$ws.Range("G251").Value = "System, dnasr281@gmail.com"
$ws.Range("G278").Value = "System, dnasr281@gmail.com"
$ws.Range("G305").Value = "System, dnasr281@gmail.com"

# --- "Not Recorded" -> "Recorded" row updates (style + values) ---
$ws.Range("A172:I172").Copy()
$ws.Range("A173:I173").PasteSpecial(-4122)
$ws.Range("G173").Value = "dnasr281@gmail.com"
$ws.Range("H173").Value = "22/23"
$ws.Range("I173").Value = "Recorded"
$ws.Range("A199:I199").Copy()
$ws.Range("A200:I200").PasteSpecial(-4122)
$ws.Range("G200").Value = "dnasr281@gmail.com"
$ws.Range("H200").Value = "26/30"
$ws.Range("I200").Value = "Recorded"
$ws.Range("A226:I226").Copy()
$ws.Range("A227:I227").PasteSpecial(-4122)
$ws.Range("G227").Value = "dnasr281@gmail.com"
$ws.Range("H227").Value = "24/25"
$ws.Range("I227").Value = "Recorded"
$ws.Range("A253:I253").Copy()
$ws.Range("A254:I254").PasteSpecial(-4122)
$ws.Range("G254").Value = "dnasr281@gmail.com"
$ws.Range("H254").Value = "23/28"
$ws.Range("I254").Value = "Recorded"
$ws.Range("A280:I280").Copy()
$ws.Range("A281:I281").PasteSpecial(-4122)
$ws.Range("G281").Value = "dnasr281@gmail.com"
$ws.Range("H281").Value = "20/26"
$ws.Range("I281").Value = "Recorded"
$ws.Range("A307:I307").Copy()
$ws.Range("A308:I308").PasteSpecial(-4122)
$ws.Range("G308").Value = "dnasr281@gmail.com"
$ws.Range("H308").Value = "22/29"
$ws.Range("I308").Value = "Recorded"
$excel.CutCopyMode = 0

